$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed coin price/volume data (scraped update)
$data = @(
    @{ Row = 2; B = 'Bitcoin'; C = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D = '30.409.66'; E = '  -0.90%  '; DForceText = $false },
    @{ Row = 3; B = 'Ethereum'; C = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D = '1.920.28'; E = '  +2.72%  '; DForceText = $false },
    @{ Row = 4; B = 'TetherUSD'; C = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D = '1.000'; E = '  +0.06%  '; DForceText = $true },
    @{ Row = 5; B = 'BNB'; C = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D = '240.99'; E = '  +1.90%  '; DForceText = $true },
    @{ Row = 6; B = 'USDC'; C = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D = '1.001'; E = '  +0.14%  '; DForceText = $true },
    @{ Row = 7; B = 'XRP'; C = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D = '0.4688'; E = '  -2.10%  '; DForceText = $true },
    @{ Row = 8; B = 'Cardano'; C = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D = '0.2859'; E = '  +0.97%  '; DForceText = $true },
    @{ Row = 9; B = 'Dogecoin'; C = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D = '0.06957'; E = '  +6.53%  '; DForceText = $true },
    @{ Row = 10; B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D = '106.87'; E = '  +12.97%  '; DForceText = $true },
    @{ Row = 11; B = 'Solana'; C = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D = '18.36'; E = '  -2.76%  '; DForceText = $true },
    @{ Row = 12; B = 'WrappedEther'; C = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D = '1.918.55'; E = '  +2.53%  '; DForceText = $false },
    @{ Row = 13; B = 'TRON'; C = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D = '0.07647'; E = '  +1.95%  '; DForceText = $true },
    @{ Row = 14; B = 'Polkadot'; C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D = '5.182'; E = '  +1.85%  '; DForceText = $true },
    @{ Row = 15; B = 'Polygon'; C = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D = '0.6574'; E = '  +1.04%  '; DForceText = $true },
    @{ Row = 16; B = 'BitcoinCash'; C = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D = '294.56'; E = '  -0.59%  '; DForceText = $true },
    @{ Row = 17; B = 'WrappedBTC'; C = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D = '30.404.46'; E = '  -0.87%  '; DForceText = $false },
    @{ Row = 18; B = 'ShibaInu'; C = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D = '0.000007650'; E = '  +2.27%  '; DForceText = $true },
    @{ Row = 19; B = 'Avalanche'; C = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D = '13.01'; E = '  +0.07%  '; DForceText = $true },
    @{ Row = 20; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '1.001'; E = '  +0.22%  '; DForceText = $true },
    @{ Row = 21; B = 'WrappedliquidstakedEther2.0'; C = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D = '2.164.00'; E = '  +2.18%  '; DForceText = $false },
    @{ Row = 22; B = 'BinanceUSD'; C = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D = '1.001'; E = '  +0.02%  '; DForceText = $true },
    @{ Row = 23; B = 'Uniswap'; C = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D = '5.221'; E = '  +0.77%  '; DForceText = $true },
    @{ Row = 24; B = 'Chainlink'; C = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D = '6.208'; E = '  +1.93%  '; DForceText = $true },
    @{ Row = 25; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '167.90'; E = '  -0.60%  '; DForceText = $true },
    @{ Row = 26; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '9.258'; E = '  +0.26%  '; DForceText = $true },
    @{ Row = 27; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '21.24'; E = '  +8.96%  '; DForceText = $true },
    @{ Row = 28; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '2.045'; E = '  +4.70%  '; DForceText = $true },
    @{ Row = 29; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D = '0.1081'; E = '  +3.37%  '; DForceText = $true },
    @{ Row = 30; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '1.364'; E = '  +0.87%  '; DForceText = $true },
    @{ Row = 31; B = 'InternetComputer(DFINITY)'; C = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D = '4.140'; E = '  +0.51%  '; DForceText = $true },
    @{ Row = 32; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '3.970'; E = '  +0.50%  '; DForceText = $true },
    @{ Row = 33; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.05055'; E = '  +1.77%  '; DForceText = $true },
    @{ Row = 34; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '0.7434'; E = '  +3.50%  '; DForceText = $true },
    @{ Row = 35; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '1.145'; E = '  -2.31%  '; DForceText = $true },
    @{ Row = 36; B = 'HuobiToken'; C = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D = '2.745'; E = '  +1.32%  '; DForceText = $true },
    @{ Row = 37; B = 'Frax'; C = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; D = '1.000'; E = '  +0.21%  '; DForceText = $true },
    @{ Row = 38; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.02019'; E = '  +4.24%  '; DForceText = $true },
    @{ Row = 39; B = 'MXToken'; C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D = '2.687'; E = '  -0.75%  '; DForceText = $true },
    @{ Row = 40; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '2.055'; E = '  +1.07%  '; DForceText = $true },
    @{ Row = 41; B = 'Quant'; C = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D = '108.43'; E = '  +1.29%  '; DForceText = $true },
    @{ Row = 42; B = 'TrustWalletToken'; C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D = '0.8724'; E = '  -1.82%  '; DForceText = $true },
    @{ Row = 43; B = 'BitcoinSV'; C = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'; D = '53.84'; E = '  +27.50%  '; DForceText = $true },
    @{ Row = 44; B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = '5.846'; E = '  +5.13%  '; DForceText = $true },
    @{ Row = 45; B = 'PaxDollar'; C = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D = '1.001'; E = '  +0.18%  '; DForceText = $true },
    @{ Row = 46; B = 'TheSandbox'; C = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D = '0.4218'; E = '  +0.64%  '; DForceText = $true },
    @{ Row = 47; B = 'Aave'; C = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D = '67.54'; E = '  +4.45%  '; DForceText = $true },
    @{ Row = 48; B = 'Aptos'; C = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D = '7.188'; E = '  -2.01%  '; DForceText = $true },
    @{ Row = 49; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '9.229'; E = '  +5.54%  '; DForceText = $true },
    @{ Row = 50; B = 'Algorand'; C = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D = '0.1209'; E = '  -1.34%  '; DForceText = $true },
    @{ Row = 51; B = 'Elrond'; C = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'; D = '34.66'; E = '  +0.04%  '; DForceText = $true }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    if ($item.DForceText) {
        # Value would otherwise be auto-coerced to a number (e.g. "1.000" -> 1);
        # the source column stores these as plain text, so force Text format first.
        $ws.Cells.Item($r, 4).NumberFormat = "@"
    }
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}

Write-Host "Updated $($data.Count) rows"